$wb = $excel.ActiveWorkbook

# Duplicate the "Croatia" sheet (template for the new market sheet) and
# place the copy right after it.
$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Copy($null, $croatia)

# The freshly created copy is now the last sheet ("Croatia (2)") - rename
# it and update its two market-specific cells.
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3205/3204"

# Match the recorded UI state: the new sheet becomes the active tab with
# B4 selected ...
$greece.Select()
$greece.Range("B4").Select()

# ... while the old Croatia sheet is left with its whole sheet selected.
$croatia.Select()
$croatia.Cells.Select()

# Re-activate Greece so it's the tab shown/selected when the file is saved.
$greece.Select()
